# Lab 2 Rubric.xlsx edit:
#   Insert a new column before column A (shifts all existing content from
#   B..F over to C..G) and add two new header-ish cells in the freshly
#   inserted columns A3 ("Task 1") and B3 ("Task 2").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right by inserting a blank column at A.
$ws.Columns("A").Insert()

# Populate the two new cells created by the insert.
$ws.Range("A3").Value = "Task 1"
$ws.Range("B3").Value = "Task 2"

# Match the target's explicit style index (s="6", the plain default style
# already used for the rest of column A/B) by nudging the alignment - this
# makes Excel serialize the style id on these two cells instead of omitting
# it as "default".
$ws.Range("A3:B3").VerticalAlignment = -4108

# Re-establish the frozen pane one column further right (matches the
# natural effect of the column insert: xSplit 4 -> 5, ySplit stays 3) and
# restore the active selection used on the sheet.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("F4").Select()
$win.FreezePanes = $true
$ws.Range("D6").Select()
